$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Predictor column (column C) text values to wrap certain
# predictors with ln(...) notation, and fix a bracket typo on the
# "Livestock AB Consumption" row, per the commit:
# "fixing marginal effects (now conditional effects) to show both parts
#  of hurdle, and contour graph for interaction"

$ws.Range("C2").Value = "ln(GDP [dollars per capita])"
$ws.Range("C3").Value = "ln(ProMed Mentions [per capita])"
$ws.Range("C4").Value = "ln(Migrant Population [per capita])"
$ws.Range("C6").Value = "ln(Tourism - Inbound [per capita])"
$ws.Range("C8").Value = "ln(AB Exports [dollars per capita])"
$ws.Range("C9").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C12").Value = "Livestock AB Consumption [kg per capita)"
$ws.Range("C13").Value = "ln(ProMed Mentions [per capita])"
$ws.Range("C14").Value = "ln(Publication Bias Index [per capita])"
$ws.Range("C16").Value = "ln(Population)"
$ws.Range("C17").Value = "ln(GDP [dollars per capita])"

$wb.Save()
